# Add new columns I ("I0") and J ("IF") to the worksheet, matching the
# style of the existing header row, and fill in the data values for rows 2-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1), copying the style used by the other headers (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (rows 2-39)
$iValues = @(5,9,7,8,7,4,6,8,5,8,8,7,7,8,6,7,6,9,8,8,7,8,7,9,8,9,8,6,6,5,4,2,5,6,8,9,8,7)

# Data values for column J (rows 2-39)
$jValues = @(5,9,7,9,7,5,6,9,6,8,8,7,8,9,7,8,6,9,8,8,9,8,8,9,9,9,8,7,6,5,5,3,6,6,9,9,8,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
